$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 4-122 (UseCode_Description, Petaluma Valley, Santa Rosa Plain,
# Sonoma Valley, Total), re-sorted descending by Total with several label/value fixes.
$data = @(
    @(4, 'Dairy w/residence', 120, 145, 30, 295),
    @(5, 'Wholesale nursery', 10, 75, 0, 85),
    @(6, 'Retail nursery', 15, 40, 15, 70),
    @(7, 'Dairy', 10, 55, 5, 70),
    @(8, 'Miscellaneous state property', 2, 40, 12, 54),
    @(9, 'Light manuftg & industrial', 0, 30, 20, 50),
    @(10, 'Auto and truck repair & maint', 0, 40, 0, 40),
    @(11, 'Religious building', 4, 28, 6, 38),
    @(12, 'Specialty shop (tires, brakes)', 0, 38, 0, 38),
    @(13, 'Store w/res unit or units', 8, 20, 10, 38),
    @(14, 'Warehousing/active', 2, 28.5, 3, 33.5),
    @(15, 'Chicken ranch w/residence', 16, 12, 4, 32),
    @(16, 'Sand & gravel, shale', 0, 30, 0, 30),
    @(17, 'Com''l use/no other category', 6, 16, 6, 28),
    @(18, 'Misc multiple use/none dominat', 4, 22, 0, 26),
    @(19, 'Light manufctrg & warehousing', 2, 24, 0, 26),
    @(20, 'Single story store', 2, 20, 2, 24),
    @(21, 'Used car lot', 0, 24, 0, 24),
    @(22, 'Miscellaneous federal property', 0, 20, 4, 24),
    @(23, 'One story office building', 4, 20, 0, 24),
    @(24, 'Municipal utility property', 8, 6, 8, 22),
    @(25, 'Single live/work unit', 0, 22, 0, 22),
    @(26, 'Winery w/vineyards', 0, 20, 0, 20),
    @(27, 'Sand and gravel, shale', 0, 0, 20, 20),
    @(28, 'Two story office building', 2, 18, 0, 20),
    @(29, 'Prop used along w/rel bldg', 2, 18, 0, 20),
    @(30, 'Veterinary hospitals', 5, 10, 5, 20),
    @(31, 'Mineral processing', 0, 10, 10, 20),
    @(32, 'Horse ranch w/residence', 0, 18, 0, 18),
    @(33, 'County park/other rec facility', 0, 6, 10, 16),
    @(34, 'Auto & truck repair & maint', 4, 0, 12, 16),
    @(35, 'Medical offices', 0, 15, 0, 15),
    @(36, 'Neighborhood shopping center', 4, 10, 0, 14),
    @(37, 'Sbe-valued utility', 2.5, 7, 3.5, 13),
    @(38, 'Other sales: trailers, mbh, rv', 0, 12, 0, 12),
    @(39, 'Warehousing yard', 0, 10, 2, 12),
    @(40, 'Privately owned park', 2, 10, 0, 12),
    @(41, 'Rest home', 2, 10, 0, 12),
    @(42, 'Club/lodge hall', 0, 8, 2, 10),
    @(43, 'Meat products', 0, 6, 4, 10),
    @(44, 'Home for handicapped (physical, mental, etc.)', 0, 10, 0, 10),
    @(45, 'Retail lumber yard', 2, 6, 2, 10),
    @(46, 'County hospital', 0, 10, 0, 10),
    @(47, 'Industr''l in no other category', 2, 8, 0, 10),
    @(48, 'Airport/private', 0, 0, 8, 8),
    @(49, 'Specialty lumber products', 0, 2, 6, 8),
    @(50, 'Truck terminal', 2, 6, 0, 8),
    @(51, 'Auto sales w/o service center', 0, 8, 0, 8),
    @(52, 'Winery with vineyards', 0, 0, 8, 8),
    @(53, 'Horse ranch w/res', 6, 0, 2, 8),
    @(54, 'Multiple stores in 1 structure', 0, 8, 0, 8),
    @(55, 'Chicken ranch', 6, 2, 0, 8),
    @(56, 'Full service station', 2, 6, 0, 8),
    @(57, 'Restaurant', 0, 6, 2, 8),
    @(58, 'Cemetery', 0, 6.5, 0.5, 7),
    @(59, 'Multiple combo/stores & office', 2, 4, 0, 6),
    @(60, 'Other food processing plants', 0, 6, 0, 6),
    @(61, 'Recreational center', 2, 2, 2, 6),
    @(62, 'Misc multiple use/no dominate', 2, 4, 0, 6),
    @(63, '18 hole public golf course', 4, 2, 0, 6),
    @(64, 'Mini-warehouse', 0.5, 4.5, 1, 6),
    @(65, 'Bulk plant', 4, 0, 2, 6),
    @(66, 'County building', 0, 6, 0, 6),
    @(67, 'Fire district', 0, 2, 4, 6),
    @(68, 'Horse ranch w/2 or more residences', 0, 6, 0, 6),
    @(69, 'Alternate use office bldgs', 0, 6, 0, 6),
    @(70, 'Industrial common area', 0, 6, 0, 6),
    @(71, 'Grocery store', 0, 4, 2, 6),
    @(72, 'Dairy w/manufactured home', 5, 0, 0, 5),
    @(73, 'Dental offices', 0, 5, 0, 5),
    @(74, 'Arcades & amusement center', 0, 4, 0, 4),
    @(75, 'State pk/other recreation fac', 2, 0, 2, 4),
    @(76, 'Cocktail lounge bar', 0, 4, 0, 4),
    @(77, 'Sfd converted to residential care facility', 0, 4, 0, 4),
    @(78, 'Service station/mini-mart', 0, 4, 0, 4),
    @(79, 'Winery', 0, 2, 2, 4),
    @(80, 'Drive-in restaurant', 2, 2, 0, 4),
    @(81, 'Feed and grain mill', 4, 0, 0, 4),
    @(82, '3-or-more story office bldg', 0, 4, 0, 4),
    @(83, 'Multiple story store', 0, 2, 2, 4),
    @(84, 'Indiv parcel/neighborhd shop ctr', 0, 4, 0, 4),
    @(85, 'Horse ranch w/2 or more res', 2, 0, 2, 4),
    @(86, 'Horse ranch', 4, 0, 0, 4),
    @(87, 'Multi-offices/residential units', 0, 2, 0, 2),
    @(88, 'Assisted care facility', 0, 2, 0, 2),
    @(89, 'Supermarket', 0, 2, 0, 2),
    @(90, 'Auto sales w/service center', 0, 2, 0, 2),
    @(91, 'Alternate use', 0, 2, 0, 2),
    @(92, 'Alternate use stores', 0, 2, 0, 2),
    @(93, 'State building', 0, 2, 0, 2),
    @(94, 'Specialty shop (tires,brakes)', 0, 0, 2, 2),
    @(95, 'Horse ranch w/manufacturedhome', 2, 0, 0, 2),
    @(96, 'Health spa or club', 2, 0, 0, 2),
    @(97, 'Chicken ranch w/manufactured home', 0, 2, 0, 2),
    @(98, 'Alternate use service stations', 0, 2, 0, 2),
    @(99, 'Alternate use prof bldgs', 0, 2, 0, 2),
    @(100, 'Sfd converted to res care fac', 0, 0, 2, 2),
    @(101, 'Live/work units', 0, 2, 0, 2),
    @(102, 'Self service sta/no repair facilities', 0, 2, 0, 2),
    @(103, 'Community shopping center', 0, 2, 0, 2),
    @(104, 'Convenience store', 0, 2, 0, 2),
    @(105, 'Country club', 0, 2, 0, 2),
    @(106, 'Lumber mill', 0, 2, 0, 2),
    @(107, 'Rural res/manufactured home', 0, 0, 2, 2),
    @(108, 'Heavy industry', 0, 2, 0, 2),
    @(109, 'Retail feed and grain sales', 0, 0, 2, 2),
    @(110, 'Other poultry ranch', 0, 2, 0, 2),
    @(111, 'Other poultry ranch w/residence', 0, 2, 0, 2),
    @(112, 'Farm or const mach sales/serv', 0, 2, 0, 2),
    @(113, 'Regional shopping center', 0, 2, 0, 2),
    @(114, 'Alternate use store/off combo', 0, 2, 0, 2),
    @(115, 'Rural residential w/misc residential imp', 0, 1.62, 0, 1.62),
    @(116, 'Utility water company', 0.5, 1, 0, 1.5),
    @(117, 'Warehousing/inactive', 0, 1, 0, 1),
    @(118, 'Mortuary/funeral home', 0, 0.5, 0, 0.5),
    @(119, 'Radio & tv broadcast site', 0, 0.5, 0, 0.5),
    @(120, 'Volunteer fire department', 0, 0, 0.5, 0.5),
    @(121, 'Mutual water company', 0, 0.5, 0, 0.5),
    @(122, 'Cable tv', 0, 0.5, 0, 0.5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Remove the 3 now-unused trailing rows (125 -> 122 data rows total).
$ws.Range("A123:E125").EntireRow.Delete()

Write-Host "Applied $($data.Count) data rows; deleted trailing rows 123-125."
